$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9.054883713581678
$ws.Range("D2").Value = 17.86437773734367
$ws.Range("E2").Value = 8.479206092280677
$ws.Range("F2").Value = 83.87726448752882
$ws.Range("G2").Value = 3.873400959496192
$ws.Range("L2").Value = 6.734303653953861
$ws.Range("M2").Value = 48.89542058806997
$ws.Range("N2").Value = 16.97149281301617

$ws.Range("C3").Value = 9.119645376577187
$ws.Range("D3").Value = 17.77270215548539
$ws.Range("E3").Value = 8.352235548181179
$ws.Range("F3").Value = 82.78789453517027
$ws.Range("G3").Value = 3.889279887184919
$ws.Range("L3").Value = 6.610293885266783
$ws.Range("M3").Value = 47.96496274182373
$ws.Range("N3").Value = 16.9895752332252

$ws.Range("C4").Value = 9.161018334428389
$ws.Range("D4").Value = 17.7256687217048
$ws.Range("E4").Value = 8.273101454849805
$ws.Range("F4").Value = 82.16067948472528
$ws.Range("G4").Value = 3.899415910356884
$ws.Range("L4").Value = 6.532729597756131
$ws.Range("M4").Value = 47.40447409730797
$ws.Range("N4").Value = 17.00359069467425

$ws.Range("C5").Value = 9.178286854543675
$ws.Range("D5").Value = 17.70879075466516
$ws.Range("E5").Value = 8.240570488195523
$ws.Range("F5").Value = 81.91560454792597
$ws.Range("G5").Value = 3.903645109149241
$ws.Range("L5").Value = 6.501047493684791
$ws.Range("M5").Value = 47.17901784300856
$ws.Range("N5").Value = 17.01004107109063

$ws.Range("C6").Value = 9.181179081982865
$ws.Range("D6").Value = 17.70612520905183
$ws.Range("E6").Value = 8.235152001589597
$ws.Range("F6").Value = 81.875546195862
$ws.Range("G6").Value = 3.904353369201953
$ws.Range("L6").Value = 6.498487811756164
$ws.Range("M6").Value = 47.14176561507856
$ws.Range("N6").Value = 17.01115696945483

$ws.Range("C7").Value = 9.161249563535549
$ws.Range("D7").Value = 17.72543188683132
$ws.Range("E7").Value = 8.272663861297206
$ws.Range("F7").Value = 82.15733168370937
$ws.Range("G7").Value = 3.899472545214959
$ws.Range("L7").Value = 6.532300051025315
$ws.Range("M7").Value = 47.40142128263845
$ws.Range("N7").Value = 17.00367468618281

$ws.Range("C8").Value = 9.076882083877431
$ws.Range("D8").Value = 17.83082444705643
$ws.Range("E8").Value = 8.435677221364422
$ws.Range("F8").Value = 83.4930367360198
$ws.Range("G8").Value = 3.87879678795086
$ws.Range("L8").Value = 6.691847700478862
$ws.Range("M8").Value = 48.57249080649248
$ws.Range("N8").Value = 16.97712705848167

$ws.Range("C9").Value = 8.924012757958955
$ws.Range("D9").Value = 18.11254367756685
$ws.Range("E9").Value = 8.745547811030928
$ws.Range("F9").Value = 86.4399036907884
$ws.Range("G9").Value = 3.841243808423195
$ws.Range("L9").Value = 6.992920893233978
$ws.Range("M9").Value = 50.94481140646425
$ws.Range("N9").Value = 16.94788181663491

$ws.Range("C10").Value = 8.819088851957279
$ws.Range("D10").Value = 18.36728979011151
$ws.Range("E10").Value = 8.966652597488363
$ws.Range("F10").Value = 88.79944704059356
$ws.Range("G10").Value = 3.815371551050378
$ws.Range("L10").Value = 7.206318067140344
$ws.Range("M10").Value = 52.72101166956437
$ws.Range("N10").Value = 16.93988704435279

$ws.Range("C11").Value = 8.772898015822376
$ws.Range("D11").Value = 18.4939325986722
$ws.Range("E11").Value = 9.065710261166323
$ws.Range("F11").Value = 89.91372515352612
$ws.Range("G11").Value = 3.803950067547483
$ws.Range("L11").Value = 7.301599614637618
$ws.Range("M11").Value = 53.53369250837379
$ws.Range("N11").Value = 16.93908486574248

$ws.Range("C12").Value = 8.755622886390892
$ws.Range("D12").Value = 18.54346327260854
$ws.Range("E12").Value = 9.102994007377657
$ws.Range("F12").Value = 90.34142873103013
$ws.Range("G12").Value = 3.799672960814882
$ws.Range("L12").Value = 7.33741503958248
$ws.Range("M12").Value = 53.8419069435699
$ws.Range("N12").Value = 16.93918034239062

$ws.Range("C13").Value = 8.75933385554198
$ws.Range("D13").Value = 18.53272555065783
$ws.Range("E13").Value = 9.094974508424167
$ws.Range("F13").Value = 90.24906144136231
$ws.Range("G13").Value = 3.800592013406983
$ws.Range("L13").Value = 7.32971346808328
$ws.Range("M13").Value = 53.77550961208858
$ws.Range("N13").Value = 16.93914216056591

$ws.Range("C14").Value = 8.771472474452796
$ws.Range("D14").Value = 18.49797582792299
$ws.Range("E14").Value = 9.068782213823964
$ws.Range("F14").Value = 89.94879808572882
$ws.Range("G14").Value = 3.803597238463091
$ws.Range("L14").Value = 7.304551533752895
$ws.Range("M14").Value = 53.55904118399153
$ws.Range("N14").Value = 16.9390847684875

$ws.Range("C15").Value = 8.778935749231071
$ws.Range("D15").Value = 18.47689634204258
$ws.Range("E15").Value = 9.052708863418459
$ws.Range("F15").Value = 89.76562310695184
$ws.Range("G15").Value = 3.805444208118792
$ws.Range("L15").Value = 7.289104325906433
$ws.Range("M15").Value = 53.42650346025057
$ws.Range("N15").Value = 16.9391013527915

$ws.Range("C16").Value = 8.822138084675748
$ws.Range("D16").Value = 18.35923219069963
$ws.Range("E16").Value = 8.960147561240619
$ws.Range("F16").Value = 88.72743825694498
$ws.Range("G16").Value = 3.816124798702742
$ws.Range("L16").Value = 7.200054408915967
$ws.Range("M16").Value = 52.66797750267038
$ws.Range("N16").Value = 16.93999580144768

$ws.Range("C17").Value = 8.849032108731741
$ws.Range("D17").Value = 18.28982399659176
$ws.Range("E17").Value = 8.902967745347818
$ws.Range("F17").Value = 88.10092962052305
$ws.Range("G17").Value = 3.822764671088424
$ws.Range("L17").Value = 7.144959815784007
$ws.Range("M17").Value = 52.20369537730707
$ws.Range("N17").Value = 16.94126469615158

$ws.Range("C18").Value = 8.864646066013963
$ws.Range("D18").Value = 18.25091331037946
$ws.Range("E18").Value = 8.869936394800957
$ws.Range("F18").Value = 87.74443718276473
$ws.Range("G18").Value = 3.826616629950235
$ws.Range("L18").Value = 7.113102386821628
$ws.Range("M18").Value = 51.93710416764057
$ws.Range("N18").Value = 16.94226186711318

$ws.Range("C19").Value = 8.869957777024357
$ws.Range("D19").Value = 18.23791162501779
$ws.Range("E19").Value = 8.858728246352307
$ws.Range("F19").Value = 87.62440173841905
$ws.Range("G19").Value = 3.827926545481512
$ws.Range("L19").Value = 7.102287329097871
$ws.Range("M19").Value = 51.84692496353812
$ws.Range("N19").Value = 16.94264565688839

$ws.Range("C20").Value = 8.846154199580528
$ws.Range("D20").Value = 18.29710777516243
$ws.Range("E20").Value = 8.909069499995372
$ws.Range("F20").Value = 88.16722412630381
$ws.Range("G20").Value = 3.822054457043415
$ws.Range("L20").Value = 7.150842225884585
$ws.Range("M20").Value = 52.25307378462453
$ws.Range("N20").Value = 16.94110201143241

$ws.Range("C21").Value = 8.767901237756718
$ws.Range("D21").Value = 18.50813974572696
$ws.Range("E21").Value = 9.076481751571903
$ws.Range("F21").Value = 90.03683768857142
$ws.Range("G21").Value = 3.802713246442342
$ws.Range("L21").Value = 7.311949467068181
$ws.Range("M21").Value = 53.62261199429641
$ws.Range("N21").Value = 16.93909085912128

$ws.Range("C22").Value = 8.71801689289957
$ws.Range("D22").Value = 18.65524599032747
$ws.Range("E22").Value = 9.184567899140408
$ws.Range("F22").Value = 91.29219810346913
$ws.Range("G22").Value = 3.790351156320212
$ws.Range("L22").Value = 7.415691523165605
$ws.Range("M22").Value = 54.52034075082859
$ws.Range("N22").Value = 16.94009940759761

$ws.Range("C23").Value = 8.744527677333211
$ws.Range("D23").Value = 18.57588375208577
$ws.Range("E23").Value = 9.127004112755122
$ws.Range("F23").Value = 90.61917013294358
$ws.Range("G23").Value = 3.796924267694087
$ws.Range("L23").Value = 7.360466504455307
$ws.Range("M23").Value = 54.04102571058205
$ws.Range("N23").Value = 16.93935143884185

$ws.Range("C24").Value = 8.847454827782938
$ws.Range("D24").Value = 18.29381168791894
$ws.Range("E24").Value = 8.906311387953416
$ws.Range("F24").Value = 88.13724085880257
$ws.Range("G24").Value = 3.822375436862647
$ws.Range("L24").Value = 7.148183357101806
$ws.Range("M24").Value = 52.23074876080928
$ws.Range("N24").Value = 16.94117472684233

$ws.Range("C25").Value = 8.964049781205198
$ws.Range("D25").Value = 18.02807818334438
$ws.Range("E25").Value = 8.662832002927173
$ws.Range("F25").Value = 85.60810524259581
$ws.Range("G25").Value = 3.851093228520501
$ws.Range("L25").Value = 6.912802812720955
$ws.Range("M25").Value = 50.29628120447003
$ws.Range("N25").Value = 16.95339294656354
